$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 3) down to the
# three new rows so the new cells share the same style index instead of
# getting a freshly-minted one.
$ws.Range("A3:G3").Copy()
$ws.Range("A4:G4").PasteSpecial(-4122)
$ws.Range("A5:G5").PasteSpecial(-4122)
$ws.Range("A6:G6").PasteSpecial(-4122)

$ws.Range("A4").Value = "arroz branco"
$ws.Range("B4").Value = "comida"
$ws.Range("C4").Value = "carboidrato"
$ws.Range("D4").Value = 130.0
$ws.Range("E4").Value = 0.4
$ws.Range("F4").Value = 2.6
$ws.Range("G4").Value = 28.2

$ws.Range("A5").Value = "leite desnatado"
$ws.Range("B5").Value = "bebida"
$ws.Range("C5").Value = "-"
$ws.Range("D5").Value = 32.0
$ws.Range("E5").Value = 0.0
$ws.Range("F5").Value = 3.1
$ws.Range("G5").Value = 49.0

$ws.Range("A6").Value = "caminhada"
$ws.Range("B6").Value = "exercícios"
$ws.Range("C6").Value = "cardio"
$ws.Range("D6").Value = 5.0
$ws.Range("E6").Value = 0.0
$ws.Range("F6").Value = 0.0
$ws.Range("G6").Value = 0.0
